# Auto-generated from changes.py -- applies cell-level numeric
# updates to the Leviathan_Profits workbook (8 sheets: ALC, ARM,
# BSM, CRP, CUL, GSM, LTW, WVR) per the target commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 3136.7693
$ws.Range("I6").Value = 309.22223
$ws.Range("J6").Value = 9498.75
$ws.Range("K6").Value = 927.66669
$ws.Range("L6").Value = 28496.25
$ws.Range("M6").Value = -815.66669
$ws.Range("N6").Value = -28720.25
$ws.Range("H17").Value = 24575.135
$ws.Range("J17").Value = 24575.135
$ws.Range("L17").Value = 73725.405
$ws.Range("N17").Value = -74061.405
$ws.Range("H19").Value = 5029.1665
$ws.Range("I19").Value = 6095.0
$ws.Range("J19").Value = 3963.3333
$ws.Range("K19").Value = 6095.0
$ws.Range("L19").Value = 3963.3333
$ws.Range("M19").Value = -5920.0
$ws.Range("N19").Value = -4313.3333
$ws.Range("H32").Value = 4000.0
$ws.Range("I32").Value = 4000.0
$ws.Range("K32").Value = 4000.0
$ws.Range("M32").Value = -3674.0
$ws.Range("H39").Value = 44.666668
$ws.Range("I39").Value = 37.75
$ws.Range("J39").Value = 100.0
$ws.Range("K39").Value = 113.25
$ws.Range("L39").Value = 300.0
$ws.Range("M39").Value = 182.75
$ws.Range("N39").Value = -892.0
$ws.Range("H51").Value = 11908620.0
$ws.Range("I51").Value = 4400.2
$ws.Range("K51").Value = 4400.2
$ws.Range("M51").Value = -3916.2
$ws.Range("H70").Value = 905.6667
$ws.Range("I70").Value = 689.6667
$ws.Range("J70").Value = 1013.6667
$ws.Range("K70").Value = 2069.0001
$ws.Range("L70").Value = 3041.0001
$ws.Range("M70").Value = -1799.0001
$ws.Range("N70").Value = -3581.0001
$ws.Range("H73").Value = 905.6667
$ws.Range("I73").Value = 689.6667
$ws.Range("J73").Value = 1013.6667
$ws.Range("K73").Value = 2069.0001
$ws.Range("L73").Value = 3041.0001
$ws.Range("M73").Value = -1133.0001
$ws.Range("N73").Value = -4913.0001
$ws.Range("H76").Value = 2612.0
$ws.Range("I76").Value = 2631.0
$ws.Range("K76").Value = 2631.0
$ws.Range("M76").Value = -2316.0
$ws.Range("H79").Value = 2612.0
$ws.Range("I79").Value = 2631.0
$ws.Range("K79").Value = 2631.0
$ws.Range("M79").Value = -1539.0
$ws.Range("H87").Value = 33354.0
$ws.Range("J87").Value = 33354.0
$ws.Range("L87").Value = 33354.0
$ws.Range("N87").Value = -35850.0
$ws.Range("H90").Value = 33354.0
$ws.Range("J90").Value = 33354.0
$ws.Range("L90").Value = 100062.0
$ws.Range("N90").Value = -112542.0
$ws.Range("H106").Value = 9050.533
$ws.Range("I106").Value = 1920.7778
$ws.Range("K106").Value = 1920.7778
$ws.Range("M106").Value = -1289.7778
$ws.Range("H107").Value = 0.0
$ws.Range("I107").Value = 0.0
$ws.Range("K107").Value = 0.0
$ws.Range("M107").ClearContents()
$ws.Range("H116").Value = 9971.429
$ws.Range("J116").Value = 10800.0
$ws.Range("L116").Value = 10800.0
$ws.Range("N116").Value = -17684.0
$ws.Range("H125").Value = 1677.75
$ws.Range("I125").Value = 1299.5
$ws.Range("J125").Value = 1719.7778
$ws.Range("K125").Value = 11695.5
$ws.Range("L125").Value = 15478.0002
$ws.Range("M125").Value = -9235.5
$ws.Range("N125").Value = -20398.0002
$ws.Range("H131").Value = 941.3158
$ws.Range("I131").Value = 941.3158
$ws.Range("K131").Value = 2823.9474
$ws.Range("M131").Value = 2216.0526
$ws.Range("H132").Value = 3830.6287
$ws.Range("I132").Value = 1350.091
$ws.Range("K132").Value = 4050.273
$ws.Range("M132").Value = -1520.273
$ws.Range("H137").Value = 35920.31
$ws.Range("I137").Value = 1445.5
$ws.Range("J137").Value = 201399.4
$ws.Range("K137").Value = 4336.5
$ws.Range("L137").Value = 604198.2
$ws.Range("M137").Value = -1786.5
$ws.Range("N137").Value = -609298.2
$ws.Range("H138").Value = 1988.8148
$ws.Range("I138").Value = 1488.4
$ws.Range("J138").Value = 2614.3333
$ws.Range("K138").Value = 4465.200000000001
$ws.Range("L138").Value = 7842.999899999999
$ws.Range("M138").Value = 674.7999999999993
$ws.Range("N138").Value = -18122.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2928.0
$ws.Range("I2").Value = 2677.6667
$ws.Range("J2").Value = 4430.0
$ws.Range("K2").Value = 2677.6667
$ws.Range("L2").Value = 4430.0
$ws.Range("M2").Value = -2564.6667
$ws.Range("N2").Value = -4656.0
$ws.Range("H32").Value = 40786.19
$ws.Range("I32").Value = 24596.186
$ws.Range("J32").Value = 118138.445
$ws.Range("K32").Value = 24596.186
$ws.Range("L32").Value = 118138.445
$ws.Range("M32").Value = -24309.186
$ws.Range("N32").Value = -118712.445
$ws.Range("H45").Value = 28171.75
$ws.Range("I45").Value = 42103.8
$ws.Range("J45").Value = 18220.285
$ws.Range("K45").Value = 42103.8
$ws.Range("L45").Value = 18220.285
$ws.Range("M45").Value = -41726.8
$ws.Range("N45").Value = -18974.285
$ws.Range("H61").Value = 2970.3572
$ws.Range("I61").Value = 1843.4445
$ws.Range("K61").Value = 1843.4445
$ws.Range("M61").Value = -1631.4445
$ws.Range("H63").Value = 1998.8462
$ws.Range("I63").Value = 1998.8462
$ws.Range("K63").Value = 1998.8462
$ws.Range("M63").Value = -1312.8462
$ws.Range("H66").Value = 1998.8462
$ws.Range("I66").Value = 1998.8462
$ws.Range("K66").Value = 9994.231
$ws.Range("M66").Value = -6562.231
$ws.Range("H97").Value = 2058.1428
$ws.Range("I97").Value = 1239.4445
$ws.Range("J97").Value = 6970.3335
$ws.Range("K97").Value = 1239.4445
$ws.Range("L97").Value = 6970.3335
$ws.Range("M97").Value = -743.4445000000001
$ws.Range("N97").Value = -7962.3335
$ws.Range("H102").Value = 1330.8572
$ws.Range("I102").Value = 1305.7188
$ws.Range("K102").Value = 1305.7188
$ws.Range("M102").Value = 316.2811999999999
$ws.Range("H116").Value = 2928.0
$ws.Range("I116").Value = 2677.6667
$ws.Range("J116").Value = 4430.0
$ws.Range("K116").Value = 2677.6667
$ws.Range("L116").Value = 4430.0
$ws.Range("M116").Value = -383.6667000000002
$ws.Range("N116").Value = -9018.0
$ws.Range("H132").Value = 2597.48
$ws.Range("I132").Value = 1708.0667
$ws.Range("J132").Value = 3931.6
$ws.Range("K132").Value = 5124.2001
$ws.Range("L132").Value = 11794.8
$ws.Range("M132").Value = -2594.2001
$ws.Range("N132").Value = -16854.8
$ws.Range("H136").Value = 2970.3572
$ws.Range("I136").Value = 1843.4445
$ws.Range("K136").Value = 5530.333500000001
$ws.Range("M136").Value = -2980.333500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2928.0
$ws.Range("I3").Value = 2677.6667
$ws.Range("J3").Value = 4430.0
$ws.Range("K3").Value = 2677.6667
$ws.Range("L3").Value = 4430.0
$ws.Range("M3").Value = -2563.6667
$ws.Range("N3").Value = -4658.0
$ws.Range("H20").Value = 4119393.0
$ws.Range("I20").Value = 6537642.0
$ws.Range("J20").Value = 8369.5
$ws.Range("K20").Value = 6537642.0
$ws.Range("L20").Value = 8369.5
$ws.Range("M20").Value = -6537395.0
$ws.Range("N20").Value = -8863.5
$ws.Range("H22").Value = 417.69232
$ws.Range("I22").Value = 402.5
$ws.Range("K22").Value = 402.5
$ws.Range("M22").Value = -229.5
$ws.Range("H35").Value = 28327.834
$ws.Range("J35").Value = 30000.0
$ws.Range("L35").Value = 30000.0
$ws.Range("N35").Value = -30620.0
$ws.Range("H82").Value = 17577.0
$ws.Range("H85").Value = 17577.0
$ws.Range("H86").Value = 1858.2069
$ws.Range("I86").Value = 1540.5834
$ws.Range("K86").Value = 1540.5834
$ws.Range("M86").Value = -417.5834
$ws.Range("H89").Value = 1858.2069
$ws.Range("I89").Value = 1540.5834
$ws.Range("K89").Value = 7702.916999999999
$ws.Range("M89").Value = -2086.916999999999
$ws.Range("H134").Value = 3231.25
$ws.Range("I134").Value = 2890.0
$ws.Range("J134").Value = 3670.0
$ws.Range("K134").Value = 8670.0
$ws.Range("L134").Value = 11010.0
$ws.Range("M134").Value = -6135.0
$ws.Range("N134").Value = -16080.0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 94249.25
$ws.Range("J20").Value = 94249.25
$ws.Range("L20").Value = 94249.25
$ws.Range("N20").Value = -94721.25
$ws.Range("H22").Value = 247.7619
$ws.Range("I22").Value = 244.28572
$ws.Range("J22").Value = 254.71428
$ws.Range("K22").Value = 244.28572
$ws.Range("L22").Value = 254.71428
$ws.Range("M22").Value = 105.71428
$ws.Range("N22").Value = -954.71428
$ws.Range("H30").Value = 94249.25
$ws.Range("J30").Value = 94249.25
$ws.Range("L30").Value = 94249.25
$ws.Range("N30").Value = -94431.25
$ws.Range("H31").Value = 2085.9
$ws.Range("I31").Value = 2209.7778
$ws.Range("J31").Value = 971.0
$ws.Range("K31").Value = 2209.7778
$ws.Range("L31").Value = 971.0
$ws.Range("M31").Value = -1914.7778
$ws.Range("N31").Value = -1561.0
$ws.Range("H34").Value = 2085.9
$ws.Range("I34").Value = 2209.7778
$ws.Range("J34").Value = 971.0
$ws.Range("K34").Value = 2209.7778
$ws.Range("L34").Value = 971.0
$ws.Range("M34").Value = -2007.7778
$ws.Range("N34").Value = -1375.0
$ws.Range("H58").Value = 1726.5883
$ws.Range("I58").Value = 1489.0
$ws.Range("J58").Value = 2498.75
$ws.Range("K58").Value = 1489.0
$ws.Range("L58").Value = 2498.75
$ws.Range("M58").Value = -1286.0
$ws.Range("N58").Value = -2904.75
$ws.Range("H68").Value = 25000.0
$ws.Range("J68").Value = 25000.0
$ws.Range("L68").Value = 25000.0
$ws.Range("N68").Value = -26498.0
$ws.Range("H71").Value = 25000.0
$ws.Range("J71").Value = 25000.0
$ws.Range("L71").Value = 75000.0
$ws.Range("N71").Value = -82488.0
$ws.Range("H99").Value = 10322.667
$ws.Range("I99").Value = 1906.0
$ws.Range("J99").Value = 14531.0
$ws.Range("K99").Value = 1906.0
$ws.Range("L99").Value = 14531.0
$ws.Range("M99").Value = -408.0
$ws.Range("N99").Value = -17527.0
$ws.Range("H105").Value = 4222.4243
$ws.Range("I105").Value = 3860.5789
$ws.Range("K105").Value = 3860.5789
$ws.Range("M105").Value = -2113.5789
$ws.Range("H107").Value = 1660.9333
$ws.Range("I107").Value = 1532.4546
$ws.Range("K107").Value = 1532.4546
$ws.Range("M107").Value = 387.5454
$ws.Range("H124").Value = 59199.0
$ws.Range("J124").Value = 59199.0
$ws.Range("L124").Value = 59199.0
$ws.Range("N124").Value = -64109.0
$ws.Range("H126").Value = 10322.667
$ws.Range("I126").Value = 1906.0
$ws.Range("J126").Value = 14531.0
$ws.Range("K126").Value = 5718.0
$ws.Range("L126").Value = 43593.0
$ws.Range("M126").Value = -3248.0
$ws.Range("N126").Value = -48533.0
$ws.Range("H128").Value = 94249.25
$ws.Range("J128").Value = 94249.25
$ws.Range("L128").Value = 94249.25
$ws.Range("N128").Value = -104209.25
$ws.Range("H134").Value = 1784.7894
$ws.Range("I134").Value = 1636.4242
$ws.Range("K134").Value = 4909.2726
$ws.Range("M134").Value = -2374.2726
$ws.Range("H136").Value = 1726.5883
$ws.Range("I136").Value = 1489.0
$ws.Range("J136").Value = 2498.75
$ws.Range("K136").Value = 4467.0
$ws.Range("L136").Value = 7496.25
$ws.Range("M136").Value = -1917.0
$ws.Range("N136").Value = -12596.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1858.1428
$ws.Range("I3").Value = 1858.1428
$ws.Range("K3").Value = 5574.428400000001
$ws.Range("M3").Value = -5462.428400000001
$ws.Range("H7").Value = 504.75
$ws.Range("I7").Value = 384.5
$ws.Range("K7").Value = 1153.5
$ws.Range("M7").Value = -1041.5
$ws.Range("H17").Value = 115.375
$ws.Range("I17").Value = 82.90909
$ws.Range("J17").Value = 186.8
$ws.Range("K17").Value = 248.72727
$ws.Range("L17").Value = 560.4000000000001
$ws.Range("M17").Value = -79.72727000000003
$ws.Range("N17").Value = -898.4000000000001
$ws.Range("H40").Value = 3351.4443
$ws.Range("I40").Value = 12.6
$ws.Range("K40").Value = 50.4
$ws.Range("M40").Value = 18.6
$ws.Range("H42").Value = 8250.0
$ws.Range("J42").Value = 8250.0
$ws.Range("L42").Value = 24750.0
$ws.Range("N42").Value = -25818.0
$ws.Range("H55").Value = 13890517.0
$ws.Range("J55").Value = 17859142.0
$ws.Range("L55").Value = 53577426.0
$ws.Range("N55").Value = -53577780.0
$ws.Range("H92").Value = 1015.8333
$ws.Range("J92").Value = 999.0
$ws.Range("L92").Value = 2997.0
$ws.Range("N92").Value = -5493.0
$ws.Range("H94").Value = 11603.308
$ws.Range("I94").Value = 4961.25
$ws.Range("J94").Value = 14555.333
$ws.Range("K94").Value = 14883.75
$ws.Range("L94").Value = 43665.999
$ws.Range("M94").Value = -14207.75
$ws.Range("N94").Value = -45017.999
$ws.Range("H97").Value = 1491.875
$ws.Range("I97").Value = 1299.0
$ws.Range("J97").Value = 1519.4286
$ws.Range("K97").Value = 3897.0
$ws.Range("L97").Value = 4558.2858
$ws.Range("M97").Value = -3401.0
$ws.Range("N97").Value = -5550.2858
$ws.Range("H133").Value = 10051.737
$ws.Range("I133").Value = 5663.8335
$ws.Range("J133").Value = 12076.923
$ws.Range("K133").Value = 16991.5005
$ws.Range("L133").Value = 36230.769
$ws.Range("M133").Value = -11931.5005
$ws.Range("N133").Value = -46350.769
$ws.Range("H141").Value = 3472.8
$ws.Range("I141").Value = 2675.4285
$ws.Range("J141").Value = 5333.3335
$ws.Range("K141").Value = 8026.2855
$ws.Range("L141").Value = 16000.0005
$ws.Range("M141").Value = -2846.2855
$ws.Range("N141").Value = -26360.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 6167000.0
$ws.Range("I11").Value = 9292429.0
$ws.Range("K11").Value = 9292429.0
$ws.Range("M11").Value = -9292290.0
$ws.Range("H33").Value = 25000.0
$ws.Range("J33").Value = 25000.0
$ws.Range("L33").Value = 25000.0
$ws.Range("N33").Value = -25504.0
$ws.Range("H40").Value = 15000.0
$ws.Range("I40").Value = 0.0
$ws.Range("J40").Value = 15000.0
$ws.Range("K40").Value = 0.0
$ws.Range("L40").Value = 15000.0
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -15302.0
$ws.Range("H49").Value = 23553.143
$ws.Range("J49").Value = 23553.143
$ws.Range("L49").Value = 23553.143
$ws.Range("N49").Value = -23921.143
$ws.Range("H52").Value = 41933.0
$ws.Range("J52").Value = 41933.0
$ws.Range("L52").Value = 41933.0
$ws.Range("N52").Value = -42451.0
$ws.Range("H105").Value = 98674.0
$ws.Range("J105").Value = 98674.0
$ws.Range("L105").Value = 98674.0
$ws.Range("N105").Value = -105662.0
$ws.Range("H107").Value = 27327.277
$ws.Range("I107").Value = 50765.777
$ws.Range("J107").Value = 3888.7778
$ws.Range("K107").Value = 50765.777
$ws.Range("L107").Value = 3888.7778
$ws.Range("M107").Value = -48845.777
$ws.Range("N107").Value = -7728.7778
$ws.Range("H113").Value = 83335816.0
$ws.Range("I113").Value = 100001780.0
$ws.Range("J113").Value = 5995.0
$ws.Range("K113").Value = 100001780.0
$ws.Range("L113").Value = 5995.0
$ws.Range("M113").Value = -99999610.0
$ws.Range("N113").Value = -10335.0
$ws.Range("H122").Value = 1650.0
$ws.Range("I122").Value = 1533.3334
$ws.Range("J122").Value = 2000.0
$ws.Range("K122").Value = 4600.0002
$ws.Range("L122").Value = 6000.0
$ws.Range("M122").Value = -2150.0002
$ws.Range("N122").Value = -10900.0
$ws.Range("H126").Value = 3173.6667
$ws.Range("J126").Value = 3510.5
$ws.Range("L126").Value = 10531.5
$ws.Range("N126").Value = -15471.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 1103.5
$ws.Range("I13").Value = 1000.0
$ws.Range("K13").Value = 1000.0
$ws.Range("M13").Value = -860.0
$ws.Range("H23").Value = 15250.0
$ws.Range("I23").Value = 15250.0
$ws.Range("K23").Value = 15250.0
$ws.Range("M23").Value = -15020.0
$ws.Range("H33").Value = 20000.0
$ws.Range("I33").Value = 20000.0
$ws.Range("K33").Value = 20000.0
$ws.Range("M33").Value = -19710.0
$ws.Range("H40").Value = 6674.0
$ws.Range("I40").Value = 7775.0
$ws.Range("K40").Value = 7775.0
$ws.Range("M40").Value = -7639.0
$ws.Range("H42").Value = 25882.0
$ws.Range("I42").Value = 19000.0
$ws.Range("J42").Value = 32764.0
$ws.Range("K42").Value = 19000.0
$ws.Range("L42").Value = 32764.0
$ws.Range("M42").Value = -18437.0
$ws.Range("N42").Value = -33890.0
$ws.Range("H49").Value = 25882.0
$ws.Range("I49").Value = 19000.0
$ws.Range("J49").Value = 32764.0
$ws.Range("K49").Value = 19000.0
$ws.Range("L49").Value = 32764.0
$ws.Range("M49").Value = -18853.0
$ws.Range("N49").Value = -33058.0
$ws.Range("H61").Value = 20434.5
$ws.Range("I61").Value = 23721.4
$ws.Range("J61").Value = 4000.0
$ws.Range("K61").Value = 23721.4
$ws.Range("L61").Value = 4000.0
$ws.Range("M61").Value = -23519.4
$ws.Range("N61").Value = -4404.0
$ws.Range("H113").Value = 20434.5
$ws.Range("I113").Value = 23721.4
$ws.Range("J113").Value = 4000.0
$ws.Range("K113").Value = 23721.4
$ws.Range("L113").Value = 4000.0
$ws.Range("M113").Value = -21551.4
$ws.Range("N113").Value = -8340.0
$ws.Range("H132").Value = 3018.2083
$ws.Range("I132").Value = 2878.1904
$ws.Range("J132").Value = 3998.3333
$ws.Range("K132").Value = 8634.5712
$ws.Range("L132").Value = 11994.9999
$ws.Range("M132").Value = -6104.5712
$ws.Range("N132").Value = -17054.9999
$ws.Range("H136").Value = 3661.0
$ws.Range("I136").Value = 2474.75
$ws.Range("J136").Value = 4188.222
$ws.Range("K136").Value = 7424.25
$ws.Range("L136").Value = 12564.666
$ws.Range("M136").Value = -4874.25
$ws.Range("N136").Value = -17664.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 17500.0
$ws.Range("I29").Value = 17500.0
$ws.Range("J29").Value = 0.0
$ws.Range("K29").Value = 17500.0
$ws.Range("L29").Value = 0.0
$ws.Range("M29").Value = -17210.0
$ws.Range("N29").ClearContents()
$ws.Range("H42").Value = 0.0
$ws.Range("J42").Value = 0.0
$ws.Range("L42").Value = 0.0
$ws.Range("N42").ClearContents()
$ws.Range("H49").Value = 13665.833
$ws.Range("I49").Value = 5375.0
$ws.Range("K49").Value = 5375.0
$ws.Range("M49").Value = -5145.0
$ws.Range("H81").Value = 3000.0833
$ws.Range("I81").Value = 3166.7778
$ws.Range("K81").Value = 6333.5556
$ws.Range("M81").Value = -5272.5556
$ws.Range("H84").Value = 3000.0833
$ws.Range("I84").Value = 3166.7778
$ws.Range("K84").Value = 31667.778
$ws.Range("M84").Value = -26363.778
$ws.Range("H96").Value = 1851.0
$ws.Range("I96").Value = 1596.0
$ws.Range("K96").Value = 1596.0
$ws.Range("M96").Value = -223.0
$ws.Range("H120").Value = 15000.0
$ws.Range("J120").Value = 15000.0
$ws.Range("L120").Value = 15000.0
$ws.Range("N120").Value = -24676.0
$ws.Range("H126").Value = 2711.1428
$ws.Range("I126").Value = 1996.0
$ws.Range("K126").Value = 5988.0
$ws.Range("M126").Value = -3518.0
$ws.Range("H132").Value = 9690.92
$ws.Range("I132").Value = 10608.048
$ws.Range("K132").Value = 31824.144
$ws.Range("M132").Value = -29294.144

